$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Alumno) values keyed by row number
$alumno = @{
    5  = 'Jon'
    7  = 'Gaizka'
    8  = 'Jon'
    9  = 'Gaizka'
    10 = 'Gaizka'
    11 = 'Gaizka'
    12 = 'Jon'
    13 = 'Gaizka'
    14 = 'Gaizka'
    15 = 'Jon'
    16 = 'Jon'
    17 = 'Gaizka'
    18 = 'Gaizka'
    19 = 'Gaizka'
    20 = 'Gaizka'
    21 = 'Jon'
}

# Column C (Tarea) values keyed by row number
$tarea = @{
    5  = 'Estilos: Base, Create y Formularios'
    7  = 'Eliminar con botón desde detalles: Proyecto, Tarea, Empleado'
    8  = 'Estilos: Detalles de proyectos, empleados y tareas'
    9  = ' Footer hecho'
    10 = 'Modelo tarea arreglado'
    11 = 'Footer con últimos empleados, herramientas, tareas, empleados (extra)'
    12 = 'Estilos: Logos footer'
    13 = 'Segunda linea de menu horizontal, botón de editar en los detalles'
    14 = 'Ultimos tres en body, estilos de crear '
    15 = 'Estilos: Ultimos empleados, proyectos, tareas y herramientas'
    16 = 'Boton eliminar y editar en "Detalles"'
    17 = 'Extra: Tareas "en proceso" del listado de tareas'
    18 = 'Extra: Herramientas con baja disponibilidad del listado de herramientas'
    19 = 'Extra: Proyectos de alto presupuesto del listado de proyectos'
    20 = 'Extra: Empleados no disponibles del listado de empleados'
    21 = 'Correccion y limpieza en codigo'
}

# Column D (Fecha) serial date values keyed by row number
$fecha = @{
    5  = 45767
    7  = 45767
    8  = 45770
    9  = 45771
    10 = 45772
    11 = 45774
    12 = 45775
    13 = 45776
    14 = 45778
    15 = 45779
    16 = 45779
    17 = 45780
    18 = 45782
    19 = 45782
    20 = 45783
    21 = 45783
}

# Most date cells use the plain built-in "d-mmm" format, but rows 18 and 19
# keep the special localized custom format ("[$-C0A]d\-mmm;@") that rows
# 13/14 used to have before this edit.
$customDateFormatRows = @(18, 19)

foreach ($r in 5..21) {
    if ($alumno.ContainsKey($r)) {
        $ws.Cells.Item($r, 2).Value = $alumno[$r]
    }
    if ($tarea.ContainsKey($r)) {
        $ws.Cells.Item($r, 3).Value = $tarea[$r]
    }
    if ($fecha.ContainsKey($r)) {
        $ws.Cells.Item($r, 4).Value = $fecha[$r]
        if ($customDateFormatRows -contains $r) {
            $ws.Cells.Item($r, 4).NumberFormat = '[$-C0A]d\-mmm;@'
        } else {
            $ws.Cells.Item($r, 4).NumberFormat = 'd-mmm'
        }
    }
}

$ws.Range("D24").Select() | Out-Null
